$wb = $excel.ActiveWorkbook

# --- Update last-updated timestamp on the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "19 Nov 2025, 09:18 AM"

# --- Insert a new leading entry into the "Top Losers" table ---
$losers = $wb.Worksheets.Item("Top Losers")

# Shift rows 23:76 down by inserting a new blank row at position 23,
# then drop the row that falls off the bottom of the (fixed-size) table.
$losers.Rows.Item(23).Insert()
$losers.Rows.Item(77).Delete()

# Populate the newly inserted row with the new stock entry
$losers.Range("A23").Value = "📉"
$losers.Range("B23").Value = "TMCV"
$losers.Range("C23").Value = -1.5685
$losers.Range("D23").Value = "N/A"
$losers.Range("E23").Value = "N/A"
